# Weekly fruit/vegetable price update for "Feria Lagunitas de Puerto Montt - Palta"
# Two new weekly records are inserted at rows 162-163 (pushing the existing
# records at old row 162 onward down by two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 162, shifting everything
# from the old row 162 down to rows 164+.
$ws.Rows("162:163").Insert()

# --- New row 162: "Primera" quality record dated 2021-11-23 (serial 44523) ---
$ws.Range("A162").Value = 4
$ws.Range("B162").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C162").Value = "Los Lagos"
$ws.Range("D162").Value = 44523
$ws.Range("E162").Value = 10
$ws.Range("F162").Value = "Fruta"
$ws.Range("G162").Value = 100106
$ws.Range("H162").Value = "Oleaginosos"
$ws.Range("I162").Value = 100106002
$ws.Range("J162").Value = "Palta"
$ws.Range("K162").Value = "Hass"
$ws.Range("L162").Value = "Primera"
$ws.Range("M162").Value = 400
$ws.Range("N162").Value = 3900
$ws.Range("O162").Value = 4000
$ws.Range("P162").Value = 3950
$ws.Range("Q162").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R162").Value = "Provincia de Quillota"
$ws.Range("S162").Value = 3950
$ws.Range("T162").Value = 1

# --- New row 163: "Segunda" quality record, same date ---
$ws.Range("A163").Value = 4
$ws.Range("B163").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C163").Value = "Los Lagos"
$ws.Range("D163").Value = 44523
$ws.Range("E163").Value = 10
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100106
$ws.Range("H163").Value = "Oleaginosos"
$ws.Range("I163").Value = 100106002
$ws.Range("J163").Value = "Palta"
$ws.Range("K163").Value = "Hass"
$ws.Range("L163").Value = "Segunda"
$ws.Range("M163").Value = 200
$ws.Range("N163").Value = 3500
$ws.Range("O163").Value = 3500
$ws.Range("P163").Value = 3500
$ws.Range("Q163").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R163").Value = "Provincia de Quillota"
$ws.Range("S163").Value = 3500
$ws.Range("T163").Value = 1
